# Apply "change image size of Product page" edit:
# - Insert a new "description" column (empty data) between "name" and "imageUrl"
#   on the product table, shifting imageUrl/price one column to the right.
# - Append 13 new product rows scraped from the Pexels-backed query.
# - Resize the table/query range from A1:D8 to A1:E23.
# - Hyperlink the final row's imageUrl cell to its own image URL (this also
#   creates/updates the built-in "Hyperlink" cell style, as Excel does).
# - Update the workbook's ExternalData_1 defined name to match the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output (4)")

$allRows = @(
    @('Womens', 'Blue Tanktop', 'https://i.ibb.co/7CQVJNm/blue-tank.png', 25),
    @('Womens', 'Floral Blouse', 'https://i.ibb.co/4W2DGKm/floral-blouse.png', 20),
    @('Womens', 'Floral Dress', 'https://i.ibb.co/KV18Ysr/floral-skirt.png', 80),
    @('Womens', 'Red Dots Dress', 'https://i.ibb.co/N3BN1bh/red-polka-dot-dress.png', 80),
    @('Womens', 'Striped Sweater', 'https://i.ibb.co/KmSkMbH/striped-sweater.png', 45),
    @('Womens', 'Yellow Track Suit', 'https://i.ibb.co/v1cvwNf/yellow-track-suit.png', 135),
    @('Womens', 'White Blouse', 'https://i.ibb.co/qBcrsJg/white-vest.png', 20),
    @('Womens', 'Teal Pink and Purple Button-up Shirt', 'https://images.pexels.com/photos/1839904/pexels-photo-1839904.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 25),
    @('Womens', 'Black Longsleeve Shirt with White Polkadots', 'https://images.pexels.com/photos/1021693/pexels-photo-1021693.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 25),
    @('Womens', 'White Shirt with Stripes ', 'https://images.pexels.com/photos/1844012/pexels-photo-1844012.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 20),
    @('Womens', 'Green Crop Top', 'https://images.pexels.com/photos/2071856/pexels-photo-2071856.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 80),
    @('Womens', 'Pink Collared Half-sleeved Top', 'https://images.pexels.com/photos/1036623/pexels-photo-1036623.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 80),
    @('Womens', 'Pink Longsleeve', 'https://images.pexels.com/photos/1735671/pexels-photo-1735671.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 45),
    @('Womens', 'Black Crew Neck T-shirt', 'https://images.pexels.com/photos/2331101/pexels-photo-2331101.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 135),
    @('Womens', 'White with blue and red striped Shirt', 'https://images.pexels.com/photos/975008/pexels-photo-975008.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 20),
    @('Womens', 'Mint Green colored Dress', 'https://images.pexels.com/photos/12189034/pexels-photo-12189034.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 25),
    @('Womens', 'White Longsleeve', 'https://images.pexels.com/photos/12189044/pexels-photo-12189044.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 20),
    @('Womens', 'Striped Dress', 'https://images.pexels.com/photos/12186931/pexels-photo-12186931.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 80),
    @('Womens', 'White and Yellow Scoop-neck Mini Dress', 'https://images.pexels.com/photos/884979/pexels-photo-884979.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 80),
    @('Womens', 'Red Dress', 'https://images.pexels.com/photos/9007318/pexels-photo-9007318.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2', 45)
)

# Header row: A=title, B=name, C=description (new, blank), D=imageUrl, E=price
$ws.Cells.Item(1,1).Value = "title"
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "description"
$ws.Cells.Item(1,4).Value = "imageUrl"
$ws.Cells.Item(1,5).Value = "price"

$r = 2
foreach ($row in $allRows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = ""
    $ws.Cells.Item($r,4).Value = $row[2]
    $ws.Cells.Item($r,5).Value = $row[3]
    $r = $r + 1
}

$lastRow = $r - 1

# Resize/extend the query table to cover the full new range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E" + $lastRow))

# Re-assert header text so the table's column metadata picks up the
# correct names for the newly extended column.
$ws.Cells.Item(1,1).Value = "title"
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "description"
$ws.Cells.Item(1,4).Value = "imageUrl"
$ws.Cells.Item(1,5).Value = "price"

# Hyperlink the last row's imageUrl cell to the image itself, matching the
# source workbook (also creates the "Hyperlink" cell style/font).
$lastImageCell = $ws.Cells.Item($lastRow, 4)
$ws.Hyperlinks.Add($lastImageCell, 'https://images.pexels.com/photos/9007318/pexels-photo-9007318.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=2')

# Update the hidden ExternalData_1 defined name so it reflects the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "output (4)!ExternalData_1") {
        $n.RefersTo = "='output (4)'!`$A`$1:`$E`$" + $lastRow
    }
}

# Move the selection the way Excel left it after the edit.
$null = $ws.Range("D9").Select()
